$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54; existing rows 54-76 shift down to 55-77.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly record. The
# categorical columns (market/region/category/etc.) are identical to the
# rest of this data subset, so copy them from row 55 (the row that used
# to be row 54 before the insert) and only set the values that actually
# differ for this record.
$ws.Cells.Item(54, 1).Value = 3
$ws.Cells.Item(54, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(54, 3).Value = "Coquimbo"
$ws.Cells.Item(54, 4).Value = 44795
$ws.Cells.Item(54, 4).NumberFormat = $ws.Cells.Item(55, 4).NumberFormat
$ws.Cells.Item(54, 5).Value = 5
$ws.Cells.Item(54, 6).Value = 100112035
$ws.Cells.Item(54, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 56
$ws.Cells.Item(54, 11).Value = 15000
$ws.Cells.Item(54, 12).Value = 15000
$ws.Cells.Item(54, 13).Value = 15000
$ws.Cells.Item(54, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(54, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(54, 16).Value = 1000
$ws.Cells.Item(54, 17).Value = 15
$ws.Cells.Item(54, 18).Value = "Hortaliza"
